# LogicComponentSequenceDiagram.pptx edit
#
# The sequence-diagram call that used to read "deletePerson(p)" is renamed
# to "deleteTask(p)" (the method being invoked on the Logic component was
# renamed). Locate the textbox on slide 1 that holds this call text and
# rename just the method-name portion, leaving the "(p)" argument list and
# all formatting (colour, size, dirty/err flags, etc.) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "deletePerson(p)") {
            $tr.Replace("deletePerson", "deleteTask") | Out-Null
        }
    }
}
